$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1875
$ws.Range("C2").Value = 0.565625
$ws.Range("P2").Value = 0.14375
$ws.Range("S2").Value = 0.103125
$ws.Range("C3").Value = 0.02162162162162162
$ws.Range("J3").Value = 0.01081081081081081
$ws.Range("P3").Value = 0.7513513513513513
$ws.Range("S3").Value = 0.2162162162162162
$ws.Range("P4").Value = 0.7818181818181819
$ws.Range("S4").Value = 0.2181818181818182
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.03686635944700461
$ws.Range("F6").Value = 0.08294930875576037
$ws.Range("J6").Value = 0.2764976958525346
$ws.Range("O6").Value = 0.0184331797235023
$ws.Range("Q6").Value = 0.1612903225806452
$ws.Range("R6").Value = 0.07834101382488479
$ws.Range("S6").Value = 0.3456221198156682
$ws.Range("B7").Value = 0.135
$ws.Range("D7").Value = 0.025
$ws.Range("E7").Value = 0.005
$ws.Range("F7").Value = 0.05
$ws.Range("J7").Value = 0.095
$ws.Range("O7").Value = 0.025
$ws.Range("Q7").Value = 0.195
$ws.Range("R7").Value = 0.08
$ws.Range("S7").Value = 0.39
$ws.Range("B8").Value = 0.09484536082474226
$ws.Range("D8").Value = 0.0288659793814433
$ws.Range("F8").Value = 0.05567010309278351
$ws.Range("J8").Value = 0.1030927835051546
$ws.Range("O8").Value = 0.01649484536082474
$ws.Range("Q8").Value = 0.1979381443298969
$ws.Range("R8").Value = 0.1030927835051546
$ws.Range("S8").Value = 0.4
$ws.Range("B9").Value = 0.1333333333333333
$ws.Range("D9").Value = 0.01428571428571429
$ws.Range("F9").Value = 0.05238095238095238
$ws.Range("J9").Value = 0.1
$ws.Range("O9").Value = 0.01904761904761905
$ws.Range("Q9").Value = 0.2238095238095238
$ws.Range("R9").Value = 0.07142857142857142
$ws.Range("S9").Value = 0.3857142857142857
$ws.Range("B10").Value = 0.1064610866372981
$ws.Range("D10").Value = 0.02422907488986784
$ws.Range("F10").Value = 0.06387665198237885
$ws.Range("J10").Value = 0.1145374449339207
$ws.Range("O10").Value = 0.0183553597650514
$ws.Range("Q10").Value = 0.2518355359765052
$ws.Range("R10").Value = 0.09985315712187959
$ws.Range("S10").Value = 0.3208516886930984
$ws.Range("G11").Value = 0.1703470031545741
$ws.Range("J11").Value = 0.08832807570977919
$ws.Range("K11").Value = 0.2018927444794953
$ws.Range("L11").Value = 0.5236593059936908
$ws.Range("S11").Value = 0.01577287066246057
$ws.Range("G12").Value = 0.6726190476190477
$ws.Range("J12").Value = 0.2619047619047619
$ws.Range("K12").Value = 0.0119047619047619
$ws.Range("L12").Value = 0.005952380952380952
$ws.Range("S12").Value = 0.04761904761904762
$ws.Range("G13").Value = 0.7647058823529411
$ws.Range("J13").Value = 0.196078431372549
$ws.Range("S13").Value = 0.0392156862745098
$ws.Range("F15").Value = 0.03097345132743363
$ws.Range("H15").Value = 0.1460176991150443
$ws.Range("I15").Value = 0.07079646017699115
$ws.Range("J15").Value = 0.3185840707964602
$ws.Range("K15").Value = 0.04424778761061947
$ws.Range("M15").Value = 0.02654867256637168
$ws.Range("O15").Value = 0.1194690265486726
$ws.Range("S15").Value = 0.2433628318584071
$ws.Range("F16").Value = 0.01339285714285714
$ws.Range("H16").Value = 0.1607142857142857
$ws.Range("I16").Value = 0.08928571428571429
$ws.Range("J16").Value = 0.3705357142857143
$ws.Range("K16").Value = 0.1116071428571429
$ws.Range("M16").Value = 0.01785714285714286
$ws.Range("O16").Value = 0.04017857142857143
$ws.Range("S16").Value = 0.1964285714285714
$ws.Range("F17").Value = 0.01612903225806452
$ws.Range("H17").Value = 0.1827956989247312
$ws.Range("I17").Value = 0.1057347670250896
$ws.Range("J17").Value = 0.4014336917562724
$ws.Range("K17").Value = 0.08960573476702509
$ws.Range("M17").Value = 0.01971326164874552
$ws.Range("O17").Value = 0.05734767025089606
$ws.Range("S17").Value = 0.1272401433691756
$ws.Range("H18").Value = 0.2173913043478261
$ws.Range("I18").Value = 0.08260869565217391
$ws.Range("J18").Value = 0.4347826086956522
$ws.Range("K18").Value = 0.08260869565217391
$ws.Range("M18").Value = 0.008695652173913044
$ws.Range("O18").Value = 0.06521739130434782
$ws.Range("S18").Value = 0.108695652173913
$ws.Range("F19").Value = 0.01453710788064269
$ws.Range("H19").Value = 0.2019892884468248
$ws.Range("I19").Value = 0.07651109410864575
$ws.Range("J19").Value = 0.3902065799540934
$ws.Range("K19").Value = 0.1101759755164499
$ws.Range("M19").Value = 0.02371843917368018
$ws.Range("N19").Value = 0.0007651109410864575
$ws.Range("O19").Value = 0.05432287681713849
$ws.Range("S19").Value = 0.1277735271614384

Write-Host "Applied all updates"
